# Apply the dated heading update and refresh all multiplication problems
$d = $word.ActiveDocument
$times = [char]215  # U+00D7 MULTIPLICATION SIGN

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-02-12 Monday" "2024-02-13 Tuesday"

$o = "61${times}98="; $n = "16${times}88="; Replace-Text $o $n
$o = "86${times}14="; $n = "48${times}79="; Replace-Text $o $n
$o = "55${times}29="; $n = "82${times}95="; Replace-Text $o $n
$o = "60${times}85="; $n = "50${times}42="; Replace-Text $o $n
$o = "87${times}35="; $n = "24${times}71="; Replace-Text $o $n
$o = "20${times}61="; $n = "70${times}96="; Replace-Text $o $n
$o = "24${times}41="; $n = "36${times}31="; Replace-Text $o $n
$o = "90${times}85="; $n = "77${times}78="; Replace-Text $o $n
$o = "63${times}93="; $n = "41${times}40="; Replace-Text $o $n
$o = "73${times}71="; $n = "37${times}33="; Replace-Text $o $n
$o = "61${times}65="; $n = "15${times}12="; Replace-Text $o $n
$o = "41${times}93="; $n = "16${times}78="; Replace-Text $o $n
$o = "58${times}57="; $n = "18${times}29="; Replace-Text $o $n
$o = "87${times}20="; $n = "98${times}20="; Replace-Text $o $n
$o = "39${times}18="; $n = "58${times}30="; Replace-Text $o $n
$o = "12${times}40="; $n = "63${times}35="; Replace-Text $o $n
$o = "85${times}26="; $n = "21${times}46="; Replace-Text $o $n
$o = "62${times}26="; $n = "87${times}80="; Replace-Text $o $n
$o = "91${times}95="; $n = "52${times}84="; Replace-Text $o $n
$o = "94${times}57="; $n = "84${times}66="; Replace-Text $o $n
$o = "16${times}48="; $n = "95${times}79="; Replace-Text $o $n
$o = "46${times}18="; $n = "30${times}48="; Replace-Text $o $n
$o = "24${times}50="; $n = "78${times}68="; Replace-Text $o $n
$o = "82${times}12="; $n = "28${times}34="; Replace-Text $o $n
$o = "98${times}13="; $n = "38${times}53="; Replace-Text $o $n
